{"js": "// Remove the leading \"\u25b2\" marker run that precedes the caption\n// \"\u8868 8-2-9 \u5be9\u6838\u7d00\u9304\", while leaving the sibling caption\n// \"\u25b2\u8868 8-2-10 \u5be9\u6838\u4eba\u54e1\" (and its own \"\u25b2\") untouched.\n\nconst results = context.document.body.search(\"\u25b2\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  const hit = results.items[i];\n  const para = hit.paragraphs.getFirstOrNullObject();\n  para.load(\"text\");\n  // eslint-disable-next-line no-await-in-loop\n  await context.sync();\n\n  if (!para.isNullObject && para.text.indexOf(\"\u8868 8-2-9 \u5be9\u6838\u7d00\u9304\") !== -1) {\n    hit.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the leading \"\u25b2\" marker character that precedes the caption\n# \"\u8868 8-2-9 \u5be9\u6838\u7d00\u9304\", while leaving the sibling caption\n# \"\u25b2\u8868 8-2-10 \u5be9\u6838\u4eba\u54e1\" (and its own \"\u25b2\") untouched.\n\n$d = $word.ActiveDocument\n$marker = [char]0x25B2   # \u25b2\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*\u8868 8-2-9 \u5be9\u6838\u7d00\u9304*\" -and $t.Length -gt 0 -and $t.Substring(0,1) -eq $marker) {\n        $start = $p.Range.Start\n        $markerRange = $d.Range($start, $start + 1)\n        $markerRange.Delete()\n    }\n}\n"}
